$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.125459666666667
$ws.Range("H2").Value = 9.376379
$ws.Range("I2").Value = 0.01427728095460815
$ws.Range("J2").Value = 0.01427728095460815
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.070065
$ws.Range("N2").Value = 3.210195
$ws.Range("O2").Value = 0.07699366399823478
$ws.Range("P2").Value = 0.07699366399823478
$ws.Range("Q2").Value = 3.344444998211667
$ws.Range("R2").Value = 30.100004983905
$ws.Range("S2").Value = 0.001099260172627496
$ws.Range("T2").Value = 0.001099260172627497

$ws.Range("G3").Value = 3.125459666666667
$ws.Range("H3").Value = 9.376379
$ws.Range("I3").Value = 0.01427728095460815
$ws.Range("J3").Value = 0.01427728095460815
$ws.Range("O3").Value = 0.08667273864337491
$ws.Range("P3").Value = 0.08667273864337491
$ws.Range("Q3").Value = 3.764883916211445
$ws.Range("R3").Value = 33.883955245903
$ws.Range("S3").Value = 0.001237451040716786
$ws.Range("T3").Value = 0.001237451040716787

$ws.Range("G4").Value = 3.125459666666667
$ws.Range("H4").Value = 9.376379
$ws.Range("I4").Value = 0.01427728095460815
$ws.Range("J4").Value = 0.01427728095460815
$ws.Range("M4").Value = 11.623441
$ws.Range("N4").Value = 34.870323
$ws.Range("O4").Value = 0.8363335973583904
$ws.Range("P4").Value = 0.8363335973583904
$ws.Range("Q4").Value = 36.32859603337966
$ws.Range("R4").Value = 326.957364300417
$ws.Range("S4").Value = 0.01194056974126387
$ws.Range("T4").Value = 0.01194056974126387

$ws.Range("I5").Value = 0.9174542056984506
$ws.Range("J5").Value = 0.9174542056984507
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.070065
$ws.Range("N5").Value = 3.210195
$ws.Range("O5").Value = 0.07699366399823478
$ws.Range("P5").Value = 0.07699366399823478
$ws.Range("Q5").Value = 214.9131294041033
$ws.Range("R5").Value = 1934.21816463693
$ws.Range("S5").Value = 0.07063816084731388
$ws.Range("T5").Value = 0.07063816084731389

$ws.Range("I6").Value = 0.9174542056984506
$ws.Range("J6").Value = 0.9174542056984507
$ws.Range("O6").Value = 0.08667273864337491
$ws.Range("P6").Value = 0.08667273864337491
$ws.Range("S6").Value = 0.07951826858776694
$ws.Range("T6").Value = 0.07951826858776694

$ws.Range("I7").Value = 0.9174542056984506
$ws.Range("J7").Value = 0.9174542056984507
$ws.Range("M7").Value = 11.623441
$ws.Range("N7").Value = 34.870323
$ws.Range("O7").Value = 0.8363335973583904
$ws.Range("P7").Value = 0.8363335973583904
$ws.Range("Q7").Value = 2334.465737832711
$ws.Range("R7").Value = 21010.1916404944
$ws.Range("S7").Value = 0.7672977762633699
$ws.Range("T7").Value = 0.7672977762633699

$ws.Range("G8").Value = 14.94475633333333
$ws.Range("H8").Value = 44.834269
$ws.Range("I8").Value = 0.06826851334694113
$ws.Range("J8").Value = 0.06826851334694113
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.070065
$ws.Range("N8").Value = 3.210195
$ws.Range("O8").Value = 0.07699366399823478
$ws.Range("P8").Value = 0.07699366399823478
$ws.Range("Q8").Value = 15.99186068582833
$ws.Range("R8").Value = 143.926746172455
$ws.Range("S8").Value = 0.005256242978293392
$ws.Range("T8").Value = 0.005256242978293392

$ws.Range("G9").Value = 14.94475633333333
$ws.Range("H9").Value = 44.834269
$ws.Range("I9").Value = 0.06826851334694113
$ws.Range("J9").Value = 0.06826851334694113
$ws.Range("O9").Value = 0.08667273864337491
$ws.Range("P9").Value = 0.08667273864337491
$ws.Range("Q9").Value = 18.00223927095922
$ws.Range("R9").Value = 162.020153438633
$ws.Range("S9").Value = 0.00591701901489118
$ws.Range("T9").Value = 0.00591701901489118

$ws.Range("G10").Value = 14.94475633333333
$ws.Range("H10").Value = 44.834269
$ws.Range("I10").Value = 0.06826851334694113
$ws.Range("J10").Value = 0.06826851334694113
$ws.Range("M10").Value = 11.623441
$ws.Range("N10").Value = 34.870323
$ws.Range("O10").Value = 0.8363335973583904
$ws.Range("P10").Value = 0.8363335973583904
$ws.Range("Q10").Value = 173.7094934998763
$ws.Range("R10").Value = 1563.385441498887
$ws.Range("S10").Value = 0.05709525135375657
$ws.Range("T10").Value = 0.05709525135375657

